# Update the "想去人数" (number of people interested) counts that changed
# between the previous data pull and the refreshed one.
#
# Sheet "展览" (sheet1) and sheet "全部类型" (sheet4) both list the same
# events, so each numeric update has to be applied on both sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# 展览 sheet
$ws1.Range("F6").Value  = 9568
$ws1.Range("F9").Value  = 1210
$ws1.Range("F10").Value = 2204
$ws1.Range("F15").Value = 453
$ws1.Range("F18").Value = 1325

# 全部类型 sheet
$ws4.Range("F7").Value  = 9568
$ws4.Range("F10").Value = 1210
$ws4.Range("F11").Value = 2204
$ws4.Range("F16").Value = 453
$ws4.Range("F19").Value = 1325
